$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.303.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.250.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.61"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.01%  "
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.96"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.05"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0793"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.599.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.253.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.748"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.219.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.43%  "
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.75"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  +5.14%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.78%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.16%  "
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.67"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.42%  "
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("E40").Value = "  +6.24%  "
$ws.Range("E41").Value = "  +6.34%  "
$ws.Range("E42").Value = "  +5.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.52"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +17.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.064.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  +3.57%  "
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +11.54%  "
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.470.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("E51").Value = "  +3.93%  "
